# Change the year in the astromap link: 2018 -> 2022, and collapse the
# heavily run-split sentence into a single plain run, matching the
# target commit.

$d = $word.ActiveDocument

$oldYear = "2018"
$newYear = "2022"
$newText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/$newYear/)."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Os mapas de*" -and $t -like "*GaNight*$oldYear*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1   # exclude the paragraph mark
        $r = $d.Range($start, $end)

        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
               '</w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'

        $r.InsertXML($xml)
    }
}
